$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (Property/Value table) ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Version bump
$ws1.Range("B3").Value = "6.0.0"

# Date bump
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value now populated
$ws1.Range("B9").Value = "Alvearie Team"

# Row 10 becomes Jurisdiction / United States of America
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# The old duplicate "Contact" row (row 11) is removed entirely, shifting
# everything below it up by one row.
$ws1.Rows.Item(11).Delete()

# --- Sheet "Elements" (structure definition table) ---
$ws2 = $wb.Worksheets.Item("Elements")

# Root extension row: Short/Definition now reflect this specific extension
$ws2.Range("K2").Value = "Employee Labor Union"
$ws2.Range("L2").Value = "Code for the labor union of the employee"
